# Auto-generated edit script for cryptos.xlsx update
# Commit: Updated symbol list on Thu Feb  2 10:50:12 UTC 2023 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, [string]$Addr, [string]$Val)
    $r = $Sheet.Range($Addr)
    # Force text storage so numeric-looking strings (prices, percentages)
    # keep their exact formatting (leading/trailing zeros, "%", "-", etc.)
    # instead of being coerced into IEEE doubles by Excels type inference.
    $r.NumberFormat = "@"
    $r.Value = $Val
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" "329.16"
Set-TextValue $ws "E2" "7.07%"
Set-TextValue $ws "D3" "40.04"
Set-TextValue $ws "E3" "8.28%"
Set-TextValue $ws "D4" "5.299"
Set-TextValue $ws "E4" "3.71%"
Set-TextValue $ws "D5" "0.08102"
Set-TextValue $ws "E5" "3.80%"
Set-TextValue $ws "E6" "5.41%"
Set-TextValue $ws "D7" "1.923"
Set-TextValue $ws "E7" "2.31%"
Set-TextValue $ws "D9" "0.9420"
Set-TextValue $ws "E9" "2.10%"
Set-TextValue $ws "E10" "25.61%"
Set-TextValue $ws "E11" "4.20%"
Set-TextValue $ws "D12" "0.09165"
Set-TextValue $ws "E12" "3.27%"
Set-TextValue $ws "D13" "0.03568"
Set-TextValue $ws "E13" "6.59%"
Set-TextValue $ws "D14" "0.09590"
Set-TextValue $ws "E14" "0.23%"
Set-TextValue $ws "D15" "0.001332"
Set-TextValue $ws "E15" "-4.11%"
Set-TextValue $ws "D16" "0.006475"
Set-TextValue $ws "E16" "11.17%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws "D17" "3.365"
Set-TextValue $ws "E17" "-1.64%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws "D18" "4.519"
Set-TextValue $ws "E18" "2.85%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue $ws "D19" "0.3511"
Set-TextValue $ws "E19" "2.60%"
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws "D20" "7.189"
Set-TextValue $ws "E20" "15.52%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue $ws "D21" "0.1332"
Set-TextValue $ws "E21" "3.59%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextValue $ws "D22" "0.2561"
Set-TextValue $ws "E22" "5.29%"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws "D23" "0.04436"
Set-TextValue $ws "E23" "2.08%"
$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue $ws "D24" "0.001222"
Set-TextValue $ws "E24" "2.29%"
$ws.Range("B25").Value = "HotbitToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue $ws "D25" "0.004339"
Set-TextValue $ws "E25" "2.04%"
Set-TextValue $ws "E26" "-7.70%"
Set-TextValue $ws "D27" "0.0003992"
Set-TextValue $ws "E27" "-0.01%"
Set-TextValue $ws "D39" "0.02480"
Set-TextValue $ws "E39" "15.28%"
Set-TextValue $ws "D40" "0.05212"
Set-TextValue $ws "E40" "3.95%"
Set-TextValue $ws "D41" "0.007643"
Set-TextValue $ws "E41" "1.51%"
Set-TextValue $ws "D42" "0.1428"
Set-TextValue $ws "E42" "6.18%"
Set-TextValue $ws "D43" "0.009105"
Set-TextValue $ws "E43" "5.06%"
Set-TextValue $ws "D44" "0.002161"
Set-TextValue $ws "E44" "6.38%"
Set-TextValue $ws "D45" "0.01095"
Set-TextValue $ws "E45" "25.31%"
Set-TextValue $ws "D46" "0.00006656"
Set-TextValue $ws "E46" "1.50%"
Set-TextValue $ws "D47" "0.00000000750"
Set-TextValue $ws "E47" "-0.01%"
Set-TextValue $ws "D48" "0.002401"
Set-TextValue $ws "E48" "139.59%"
Set-TextValue $ws "D50" "0.00002101"
Set-TextValue $ws "E50" "-0.01%"
Set-TextValue $ws "D51" "0.0002001"
Set-TextValue $ws "E51" "-0.01%"
